$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert two new rows (one at row 5, one at what becomes row 7) ---
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(7).Insert()

# --- Widen column B from 40 to 41 characters ---
$ws.Columns.Item(2).ColumnWidth = 40.15

# --- Rewrite all data rows (2-15) with the final values ---
$timestamp = "2025-09-04 13:10:28"

$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【注目】AI音声認識を活用した福祉相談支援システム開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5386901"
$ws.Range("G2").Value = 398
$ws.Range("H2").Value = "🔥AI,Ai ◆開発,システム開発"

$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【急募】LINEで買取査定のAIシステム構築をお手伝いください!"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5386178"
$ws.Range("G3").Value = 313
$ws.Range("H3").Value = "🔥AI,Ai"

$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "日本株・米国株ランキングメール自動配信システムの作成依頼。Pythonなど。"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5386223"
$ws.Range("G4").Value = 205
$ws.Range("H4").Value = "🔥Python"

$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "1688アリババの商品情報の抽出のスクレイピングの開発 exe形式の自動ツール"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5387065"
$ws.Range("G5").Value = 173
$ws.Range("H5").Value = "◆ツール,開発"

$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "【iPhoneアプリ開発】マインドを高めるMy routine管理アプリ"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5386904"
$ws.Range("G6").Value = 123
$ws.Range("H6").Value = "◆開発 ◇アプリ"

$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "[MVP開発] 公式LINE向け英単語問題配信システム開発"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5387024"
$ws.Range("G7").Value = 118
$ws.Range("H7").Value = "◆開発,システム開発"

$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "初回 スクレイピング Aliexpressの商品情報の抽出ツール"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5372687"
$ws.Range("G8").Value = 108
$ws.Range("H8").Value = "◆ツール,スクレイピング"

$ws.Range("A9").Value = $timestamp
$ws.Range("B9").Value = "【ペットのアバター化】LumiGOプロトタイプ開発の依頼"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5382213"
$ws.Range("G9").Value = 68
$ws.Range("H9").Value = "◆開発"

$ws.Range("A10").Value = $timestamp
$ws.Range("B10").Value = "オンラインスロットのスクレイピングソフトの制作"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5386440"
$ws.Range("G10").Value = 48
$ws.Range("H10").Value = "◆スクレイピング"

$ws.Range("A11").Value = $timestamp
$ws.Range("B11").Value = "【急募】RUBYからPHPへのリプレース仕様書作成依頼"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5386592"
$ws.Range("G11").Value = 40
$ws.Range("H11").Value = "○PHP"

$ws.Range("A12").Value = $timestamp
$ws.Range("B12").Value = "【緊急】運営しているサイトに表示される詐欺広告の削除方法を教えてください"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5386516"
$ws.Range("G12").Value = 30
$ws.Range("H12").Value = "◇サイト"

$ws.Range("A13").Value = $timestamp
$ws.Range("B13").Value = "限定公開 PR 限定公開の仕事"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5385681"
$ws.Range("G13").Value = 25
$ws.Range("H13").ClearContents()

$ws.Range("A14").Value = $timestamp
$ws.Range("B14").Value = "注目 PR 超初級・SE育成の技術研修 サブ講師"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5385021"
$ws.Range("G14").Value = 25
$ws.Range("H14").ClearContents()

$ws.Range("A15").Value = $timestamp
$ws.Range("B15").Value = "限定公開 限定公開の仕事"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5386235"
$ws.Range("G15").Value = 18
$ws.Range("H15").ClearContents()

# --- Rebuild hyperlinks for F2:F15 from scratch (row insert does not shift them) ---
$ws.Range("A1:H15").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5386901")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5386178")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5386223")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5387065")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5386904")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5387024")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5372687")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5382213")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5386440")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5386592")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5386516")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5385681")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5385021")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5386235")
